$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "Property1" to "DataNode"
$ws.Name = "DataNode"

# Update selection on the sheet (cursor moved from H16 to D22)
$ws.Range("D22").Select()
